$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 85.36364133333332
$ws.Cells.Item(2, 8).Value = 256.090924
$ws.Cells.Item(2, 9).Value = 0.832590152283795
$ws.Cells.Item(2, 10).Value = 0.8325901522837948
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 112.513392
$ws.Cells.Item(2, 14).Value = 337.540176
$ws.Cells.Item(2, 15).Value = 0.3275312977368564
$ws.Cells.Item(2, 16).Value = 0.3275312977368564
$ws.Cells.Item(2, 17).Value = 9604.552839884735
$ws.Cells.Item(2, 18).Value = 86440.97555896261
$ws.Cells.Item(2, 19).Value = 0.2726993330604382
$ws.Cells.Item(2, 20).Value = 0.2726993330604382

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 85.36364133333332
$ws.Cells.Item(3, 8).Value = 256.090924
$ws.Cells.Item(3, 9).Value = 0.832590152283795
$ws.Cells.Item(3, 10).Value = 0.8325901522837948
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 106.314466
$ws.Cells.Item(3, 14).Value = 318.943398
$ws.Cells.Item(3, 15).Value = 0.3094859589441663
$ws.Cells.Item(3, 16).Value = 0.3094859589441664
$ws.Cells.Item(3, 17).Value = 9075.389944168859
$ws.Cells.Item(3, 18).Value = 81678.50949751974
$ws.Cells.Item(3, 19).Value = 0.2576749616870198
$ws.Cells.Item(3, 20).Value = 0.2576749616870198

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 85.36364133333332
$ws.Cells.Item(4, 8).Value = 256.090924
$ws.Cells.Item(4, 9).Value = 0.832590152283795
$ws.Cells.Item(4, 10).Value = 0.8325901522837948
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 124.6916553333333
$ws.Cells.Item(4, 14).Value = 374.074966
$ws.Cells.Item(4, 15).Value = 0.3629827433189773
$ws.Cells.Item(4, 16).Value = 0.3629827433189773
$ws.Cells.Item(4, 17).Value = 10644.13374313429
$ws.Cells.Item(4, 18).Value = 95797.20368820858
$ws.Cells.Item(4, 19).Value = 0.302215857536337
$ws.Cells.Item(4, 20).Value = 0.302215857536337

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 11.230072
$ws.Cells.Item(5, 8).Value = 33.690216
$ws.Cells.Item(5, 9).Value = 0.1095319647872954
$ws.Cells.Item(5, 10).Value = 0.1095319647872954
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 112.513392
$ws.Cells.Item(5, 14).Value = 337.540176
$ws.Cells.Item(5, 15).Value = 0.3275312977368564
$ws.Cells.Item(5, 16).Value = 0.3275312977368564
$ws.Cells.Item(5, 17).Value = 1263.533493124224
$ws.Cells.Item(5, 18).Value = 11371.80143811801
$ws.Cells.Item(5, 19).Value = 0.03587514657045052
$ws.Cells.Item(5, 20).Value = 0.03587514657045052

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 11.230072
$ws.Cells.Item(6, 8).Value = 33.690216
$ws.Cells.Item(6, 9).Value = 0.1095319647872954
$ws.Cells.Item(6, 10).Value = 0.1095319647872954
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 106.314466
$ws.Cells.Item(6, 14).Value = 318.943398
$ws.Cells.Item(6, 15).Value = 0.3094859589441663
$ws.Cells.Item(6, 16).Value = 0.3094859589441664
$ws.Cells.Item(6, 17).Value = 1193.919107821552
$ws.Cells.Item(6, 18).Value = 10745.27197039397
$ws.Cells.Item(6, 19).Value = 0.03389860515723478
$ws.Cells.Item(6, 20).Value = 0.03389860515723478

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 11.230072
$ws.Cells.Item(7, 8).Value = 33.690216
$ws.Cells.Item(7, 9).Value = 0.1095319647872954
$ws.Cells.Item(7, 10).Value = 0.1095319647872954
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 124.6916553333333
$ws.Cells.Item(7, 14).Value = 374.074966
$ws.Cells.Item(7, 15).Value = 0.3629827433189773
$ws.Cells.Item(7, 16).Value = 0.3629827433189773
$ws.Cells.Item(7, 17).Value = 1400.296267192517
$ws.Cells.Item(7, 18).Value = 12602.66640473266
$ws.Cells.Item(7, 19).Value = 0.03975821305961011
$ws.Cells.Item(7, 20).Value = 0.03975821305961011

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 5.934092333333333
$ws.Cells.Item(8, 8).Value = 17.802277
$ws.Cells.Item(8, 9).Value = 0.05787788292890966
$ws.Cells.Item(8, 10).Value = 0.05787788292890966
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 112.513392
$ws.Cells.Item(8, 14).Value = 337.540176
$ws.Cells.Item(8, 15).Value = 0.3275312977368564
$ws.Cells.Item(8, 16).Value = 0.3275312977368564
$ws.Cells.Item(8, 17).Value = 667.6648568645279
$ws.Cells.Item(8, 18).Value = 6008.983711780752
$ws.Cells.Item(8, 19).Value = 0.01895681810596763
$ws.Cells.Item(8, 20).Value = 0.01895681810596762

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 5.934092333333333
$ws.Cells.Item(9, 8).Value = 17.802277
$ws.Cells.Item(9, 9).Value = 0.05787788292890966
$ws.Cells.Item(9, 10).Value = 0.05787788292890966
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 106.314466
$ws.Cells.Item(9, 14).Value = 318.943398
$ws.Cells.Item(9, 15).Value = 0.3094859589441663
$ws.Cells.Item(9, 16).Value = 0.3094859589441664
$ws.Cells.Item(9, 17).Value = 630.8798576130273
$ws.Cells.Item(9, 18).Value = 5677.918718517246
$ws.Cells.Item(9, 19).Value = 0.0179123920999118
$ws.Cells.Item(9, 20).Value = 0.0179123920999118

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 5.934092333333333
$ws.Cells.Item(10, 8).Value = 17.802277
$ws.Cells.Item(10, 9).Value = 0.05787788292890966
$ws.Cells.Item(10, 10).Value = 0.05787788292890966
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 124.6916553333333
$ws.Cells.Item(10, 14).Value = 374.074966
$ws.Cells.Item(10, 15).Value = 0.3629827433189773
$ws.Cells.Item(10, 16).Value = 0.3629827433189773
$ws.Cells.Item(10, 17).Value = 739.9317959441759
$ws.Cells.Item(10, 18).Value = 6659.386163497582
$ws.Cells.Item(10, 19).Value = 0.02100867272303023
$ws.Cells.Item(10, 20).Value = 0.02100867272303024
